# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.250.01"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "1.673.26"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").Value = "'218.07"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("D6").Value = "'0.5137"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("D8").Value = "'0.2662"
$ws.Range("E8").Value = "  +5.41%  "
$ws.Range("D9").Value = "'0.06399"
$ws.Range("E9").Value = "  +4.74%  "
$ws.Range("D10").Value = "'21.58"
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("D11").Value = "'0.07394"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").Value = "1.674.60"
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("D14").Value = "'0.5833"
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("D15").Value = "1.901.81"
$ws.Range("E15").Value = "  +1.61%  "
$ws.Range("D16").Value = "'0.000008731"
$ws.Range("E16").Value = "  +8.78%  "
$ws.Range("D17").Value = "'64.85"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").Value = "26.319.09"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("D19").Value = "'4.963"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "'1.007"
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").Value = "'10.86"
$ws.Range("E21").Value = "  +3.03%  "
$ws.Range("D22").Value = "'189.70"
$ws.Range("E22").Value = "  +4.18%  "
$ws.Range("D23").Value = "'6.224"
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("D24").Value = "'1.008"
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("D25").Value = "'144.49"
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("D26").Value = "'7.636"
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("D27").Value = "'0.1185"
$ws.Range("E27").Value = "  +4.10%  "
$ws.Range("D28").Value = "'15.65"
$ws.Range("E28").Value = "  +4.50%  "
$ws.Range("D29").Value = "'0.05922"
$ws.Range("E29").Value = "  +3.46%  "
$ws.Range("D30").Value = "'1.283"
$ws.Range("E30").Value = "  -3.17%  "
$ws.Range("D31").Value = "'1.321"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'3.526"
$ws.Range("E32").Value = "  +2.96%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'3.530"
$ws.Range("E33").Value = "  +3.85%  "
$ws.Range("E34").Value = "  +4.41%  "
$ws.Range("E35").Value = "  +3.61%  "
$ws.Range("D36").Value = "'0.6026"
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("D37").Value = "'2.370"
$ws.Range("E37").Value = "  -2.71%  "
$ws.Range("D38").Value = "'2.651"
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("D39").Value = "'0.01621"
$ws.Range("E39").Value = "  +3.76%  "
$ws.Range("D40").Value = "'6.080"
$ws.Range("E40").Value = "  +6.59%  "
$ws.Range("D41").Value = "1.079.36"
$ws.Range("E41").Value = "  +1.41%  "
$ws.Range("D42").Value = "'0.8706"
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("D44").Value = "'99.89"
$ws.Range("E44").Value = "  +4.44%  "
$ws.Range("D45").Value = "1.823.14"
$ws.Range("E45").Value = "  +2.23%  "
$ws.Range("D46").Value = "'0.00000000115"
$ws.Range("E46").Value = "  +5.27%  "
$ws.Range("D47").Value = "'56.11"
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("D49").Value = "'8.096"
$ws.Range("E49").Value = "  +3.60%  "
$ws.Range("D50").Value = "'0.4303"
$ws.Range("E50").Value = "  -1.84%  "
$ws.Range("D51").Value = "'0.05207"
$ws.Range("E51").Value = "  -0.01%  "
